$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13.80322018070584
$ws.Range("D2").Value = 9.963430740211422
$ws.Range("E2").Value = 14.22233482818651
$ws.Range("F2").Value = 49.57049264396144
$ws.Range("G2").Value = 3.762739942625077
$ws.Range("K2").Value = 26.66393156544392
$ws.Range("L2").Value = 9.53081723043509
$ws.Range("N2").Value = 20.92828258803572
$ws.Range("C3").Value = 13.72902614733063
$ws.Range("D3").Value = 9.997839943315011
$ws.Range("E3").Value = 14.18370385081651
$ws.Range("F3").Value = 49.30508043983612
$ws.Range("G3").Value = 3.768115448607392
$ws.Range("K3").Value = 26.29014404638485
$ws.Range("L3").Value = 9.543643255067462
$ws.Range("N3").Value = 20.99936157231682
$ws.Range("C4").Value = 13.68719552942605
$ws.Range("D4").Value = 10.0206737149923
$ws.Range("E4").Value = 14.1636433364695
$ws.Range("F4").Value = 49.15657148312435
$ws.Range("G4").Value = 3.771579468851736
$ws.Range("K4").Value = 26.0666641235117
$ws.Range("L4").Value = 9.553296380009092
$ws.Range("N4").Value = 21.04507178049883
$ws.Range("C5").Value = 13.67109563928436
$ws.Range("D5").Value = 10.03040605334168
$ws.Range("E5").Value = 14.15639069707033
$ws.Range("F5").Value = 49.09971405037496
$ws.Range("G5").Value = 3.773032379505713
$ws.Range("K5").Value = 25.97721861803707
$ws.Range("L5").Value = 9.557676068592997
$ws.Range("N5").Value = 21.06421941219732
$ws.Range("C6").Value = 13.66847968277914
$ws.Range("D6").Value = 10.03204786534108
$ws.Range("E6").Value = 14.15524214094189
$ws.Range("F6").Value = 49.09049468908292
$ws.Range("G6").Value = 3.773276134035017
$ws.Range("K6").Value = 25.96246747560656
$ws.Range("L6").Value = 9.558430213689977
$ws.Range("N6").Value = 21.06743030433119
$ws.Range("C7").Value = 13.68697455680461
$ws.Range("D7").Value = 10.02080324041046
$ws.Range("E7").Value = 14.16354178890814
$ws.Range("F7").Value = 49.15578982970774
$ws.Range("G7").Value = 3.771598895852334
$ws.Range("K7").Value = 26.06545111199051
$ws.Range("L7").Value = 9.553353641842016
$ws.Range("N7").Value = 21.04532790459396
$ws.Range("C8").Value = 13.77687131985278
$ws.Range("D8").Value = 9.974939844503343
$ws.Range("E8").Value = 14.20825696676342
$ws.Range("F8").Value = 49.47599595782126
$ws.Range("G8").Value = 3.764559614283276
$ws.Range("K8").Value = 26.53386502221132
$ws.Range("L8").Value = 9.534869999207684
$ws.Range("N8").Value = 20.95236181518629
$ws.Range("C9").Value = 13.98221465500236
$ws.Range("D9").Value = 9.898618866770809
$ws.Range("E9").Value = 14.32484300587248
$ws.Range("F9").Value = 50.21718008034706
$ws.Range("G9").Value = 3.75204342538382
$ws.Range("K9").Value = 27.49531771502282
$ws.Range("L9").Value = 9.512777681762879
$ws.Range("N9").Value = 20.78643667288109
$ws.Range("C10").Value = 14.15003287490983
$ws.Range("D10").Value = 9.85095654472517
$ws.Range("E10").Value = 14.4278788112044
$ws.Range("F10").Value = 50.82847086955582
$ws.Range("G10").Value = 3.743620260509267
$ws.Range("K10").Value = 28.22112532213283
$ws.Range("L10").Value = 9.50523669349989
$ws.Range("N10").Value = 20.6744842531921
$ws.Range("C11").Value = 14.22986987016353
$ws.Range("D11").Value = 9.831124652246173
$ws.Range("E11").Value = 14.47845596500485
$ws.Range("F11").Value = 51.12045190046405
$ws.Range("G11").Value = 3.739953363194101
$ws.Range("K11").Value = 28.55409920636929
$ws.Range("L11").Value = 9.503704603747378
$ws.Range("N11").Value = 20.62570889824559
$ws.Range("C12").Value = 14.26058702901937
$ws.Range("D12").Value = 9.823882994566373
$ws.Range("E12").Value = 14.49813359801814
$ws.Range("F12").Value = 51.23295802957003
$ws.Range("G12").Value = 3.738588296109392
$ws.Range("K12").Value = 28.68047223652705
$ws.Range("L12").Value = 9.503398173332926
$ws.Range("N12").Value = 20.60754815001287
$ws.Range("C13").Value = 14.25395028333904
$ws.Range("D13").Value = 9.825430648445947
$ws.Range("E13").Value = 14.49387244285385
$ws.Range("F13").Value = 51.20864257060622
$ws.Range("G13").Value = 3.738881245373675
$ws.Range("K13").Value = 28.65324483255714
$ws.Range("L13").Value = 9.50345198294958
$ws.Range("N13").Value = 20.61144563130953
$ws.Range("C14").Value = 14.23238738362169
$ws.Range("D14").Value = 9.830523487598027
$ws.Range("E14").Value = 14.48006437789491
$ws.Range("F14").Value = 51.12966933981807
$ws.Range("G14").Value = 3.739840588311798
$ws.Range("K14").Value = 28.56449099525033
$ws.Range("L14").Value = 9.503673903398823
$ws.Range("N14").Value = 20.62420860164287
$ws.Range("C15").Value = 14.21924204452826
$ws.Range("D15").Value = 9.833677999886817
$ws.Range("E15").Value = 14.47167468753699
$ws.Range("F15").Value = 51.08154665450284
$ws.Range("G15").Value = 3.740431269059409
$ws.Range("K15").Value = 28.5101600599677
$ws.Range("L15").Value = 9.503845505814031
$ws.Range("N15").Value = 20.63206658340817
$ws.Range("C16").Value = 14.14488416092639
$ws.Range("D16").Value = 9.852290020310335
$ws.Range("E16").Value = 14.42464740377152
$ws.Range("F16").Value = 50.80966394913637
$ws.Range("G16").Value = 3.743863201115876
$ws.Range("K16").Value = 28.19941076450156
$ws.Range("L16").Value = 9.505375085815421
$ws.Range("N16").Value = 20.67771513974667
$ws.Range("C17").Value = 14.10015080833689
$ws.Range("D17").Value = 9.86418327321854
$ws.Range("E17").Value = 14.39674159220981
$ws.Range("F17").Value = 50.64639214867022
$ws.Range("G17").Value = 3.746010660552404
$ws.Range("K17").Value = 28.00940752722104
$ws.Range("L17").Value = 9.506800215920439
$ws.Range("N17").Value = 20.7062702436834
$ws.Range("C18").Value = 14.07475154218981
$ws.Range("D18").Value = 9.871197945443184
$ws.Range("E18").Value = 14.38104021935738
$ws.Range("F18").Value = 50.55379490462071
$ws.Range("G18").Value = 3.747261349996491
$ws.Range("K18").Value = 27.90039331179093
$ws.Range("L18").Value = 9.507798554935253
$ws.Range("N18").Value = 20.72289697685333
$ws.Range("C19").Value = 14.06620902834271
$ws.Range("D19").Value = 9.873602808458367
$ws.Range("E19").Value = 14.37578421928244
$ws.Range("F19").Value = 50.52267022072557
$ws.Range("G19").Value = 3.747687484812619
$ws.Range("K19").Value = 27.86353313456089
$ws.Range("L19").Value = 9.508167232543151
$ws.Range("N19").Value = 20.72856130190206
$ws.Range("C20").Value = 14.10487870578927
$ws.Range("D20").Value = 9.86289919311665
$ws.Range("E20").Value = 14.3996761140755
$ws.Range("F20").Value = 50.66363730395975
$ws.Range("G20").Value = 3.745780454208641
$ws.Range("K20").Value = 28.02960656228396
$ws.Range("L20").Value = 9.506630013829033
$ws.Range("N20").Value = 20.703209537362
$ws.Range("C21").Value = 14.23870792616153
$ws.Range("D21").Value = 9.829020298982639
$ws.Range("E21").Value = 14.48410595424259
$ws.Range("F21").Value = 51.15281355923464
$ws.Range("G21").Value = 3.739558169624834
$ws.Range("K21").Value = 28.59055340782293
$ws.Range("L21").Value = 9.503601285225423
$ws.Range("N21").Value = 20.62045141072534
$ws.Range("C22").Value = 14.32898876923323
$ws.Range("D22").Value = 9.808443125853552
$ws.Range("E22").Value = 14.54234260821584
$ws.Range("F22").Value = 51.48379128947094
$ws.Range("G22").Value = 3.73562847396561
$ws.Range("K22").Value = 28.95876086758516
$ws.Range("L22").Value = 9.503217677990467
$ws.Range("N22").Value = 20.56816772523521
$ws.Range("C23").Value = 14.2805527552102
$ws.Range("D23").Value = 9.819281632448922
$ws.Range("E23").Value = 14.51098367346628
$ws.Range("F23").Value = 51.30613162880733
$ws.Range("G23").Value = 3.737713362849974
$ws.Range("K23").Value = 28.76213439297939
$ws.Range("L23").Value = 9.503276171108393
$ws.Range("N23").Value = 20.59590750358236
$ws.Range("C24").Value = 14.102740229831
$ws.Range("D24").Value = 9.863479174325736
$ws.Range("E24").Value = 14.39834835075636
$ws.Range("F24").Value = 50.65583680934825
$ws.Range("G24").Value = 3.745884480375864
$ws.Range("K24").Value = 28.0204738805128
$ws.Range("L24").Value = 9.506706404668233
$ws.Range("N24").Value = 20.70459262850644
$ws.Range("C25").Value = 13.92362641806514
$ws.Range("D25").Value = 9.917796825550628
$ws.Range("E25").Value = 14.29023020232255
$ws.Range("F25").Value = 50.00476045402417
$ws.Range("G25").Value = 3.755292818234081
$ws.Range("K25").Value = 27.23131811722122
$ws.Range("L25").Value = 9.51723224210189
$ws.Range("N25").Value = 20.82957375051956
